# Update the "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect refreshed scrape data, per commit:
#   "Update gh-pages to output generated at 456a3b4"

$wb = $excel.ActiveWorkbook

# Row -> new F-column value for the "展览" sheet
$sheet1Updates = @{
    2  = 262
    3  = 1095
    5  = 464
    6  = 90
    7  = 581
    8  = 82
    9  = 6912
    10 = 170
    12 = 140
    15 = 1124
    16 = 16407
    17 = 7
    18 = 1612
    20 = 346
    21 = 192
    23 = 11466
    25 = 1102
    26 = 4515
    27 = 371
    29 = 56
    33 = 5215
}

# Row -> new F-column value for the "全部类型" sheet
$sheet4Updates = @{
    2  = 262
    3  = 1095
    5  = 464
    6  = 90
    7  = 581
    9  = 82
    10 = 6912
    11 = 170
    13 = 140
    17 = 1124
    18 = 16407
    19 = 7
    20 = 1612
    22 = 346
    23 = 192
    27 = 11466
    29 = 1102
    30 = 4515
    31 = 371
    33 = 56
    37 = 5215
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
